$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove all existing hyperlinks so stale rId/targets do not linger
$ws.Range("F2").Hyperlinks.Delete()

# Row 2
$ws.Range("A2").Value = '2026-02-09 07:01:52'
$ws.Range("B2").Value = '【未来予測】パラレルワールドAIシステム開発の依頼'
$ws.Range("C2").Value = 'システム開発'
$ws.Range("D2").Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Range("E2").Value = '期限情報なし'
$ws.Range("F2").Value = 'https://www.lancers.jp/work/detail/5488301'
$ws.Range("G2").Value = 403
$ws.Range("H2").Value = '🔥AI,Ai ◆開発,システム開発'

# Row 3
$ws.Range("A3").Value = '2026-02-09 07:01:52'
$ws.Range("B3").Value = '【急募】AIロボット・エージェント動作生成&販売プラットフォーム開発'
$ws.Range("C3").Value = 'システム開発'
$ws.Range("D3").Value = '200,000 円 ~ 300,000 円 / 固定'
$ws.Range("E3").Value = '期限情報なし'
$ws.Range("F3").Value = 'https://www.lancers.jp/work/detail/5488299'
$ws.Range("G3").Value = 368
$ws.Range("H3").Value = '🔥AI,Ai ◆開発'

# Row 4
$ws.Range("A4").Value = '2026-02-09 07:01:52'
$ws.Range("B4").Value = '【急募】パーソナルAI開発プロジェクトの依頼'
$ws.Range("C4").Value = 'システム開発'
$ws.Range("D4").Value = '200,000 円 ~ 300,000 円 / 固定'
$ws.Range("E4").Value = '期限情報なし'
$ws.Range("F4").Value = 'https://www.lancers.jp/work/detail/5488286'
$ws.Range("G4").Value = 368
$ws.Range("H4").Value = '🔥AI,Ai ◆開発'

# Row 5
$ws.Range("A5").Value = '2026-02-09 07:01:52'
$ws.Range("B5").Value = '【急募】AIシミュレーション相性チェックサービス開発者募集'
$ws.Range("C5").Value = 'システム開発'
$ws.Range("D5").Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Range("E5").Value = '期限情報なし'
$ws.Range("F5").Value = 'https://www.lancers.jp/work/detail/5488266'
$ws.Range("G5").Value = 368
$ws.Range("H5").Value = '🔥AI,Ai ◆開発'

# Row 6
$ws.Range("A6").Value = '2026-02-09 07:01:52'
$ws.Range("B6").Value = '【無在庫ツール開発】KeepaAPIとbaseAPIを活用したシステム構築'
$ws.Range("C6").Value = 'システム開発'
$ws.Range("D6").Value = '1,000 ~ 5,000 円 / 固定'
$ws.Range("E6").Value = '期限情報なし'
$ws.Range("F6").Value = 'https://www.lancers.jp/work/detail/5488392'
$ws.Range("G6").Value = 320
$ws.Range("H6").Value = '🔥API ◆ツール,開発'

# Row 7
$ws.Range("A7").Value = '2026-02-09 07:01:52'
$ws.Range("B7").Value = '※急募:Flutterによる業務アプリの開発(+next.js)'
$ws.Range("C7").Value = 'システム開発'
$ws.Range("D7").Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Range("E7").Value = '期限情報なし'
$ws.Range("F7").Value = 'https://www.lancers.jp/work/detail/5488271'
$ws.Range("G7").Value = 225
$ws.Range("H7").Value = '🔥Next.js ◆開発 ◇アプリ'

# Row 8
$ws.Range("A8").Value = '2026-02-09 07:01:52'
$ws.Range("B8").Value = '【Zapier保守・運用サポート】既存フローの管理・調整をお任せできる方募集(時給1,200円程度)'
$ws.Range("C8").Value = 'システム開発'
$ws.Range("D8").Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Range("E8").Value = '期限情報なし'
$ws.Range("F8").Value = 'https://www.lancers.jp/work/detail/5488168'
$ws.Range("G8").Value = 213
$ws.Range("H8").Value = '🔥API ◇管理'

# Row 9
$ws.Range("A9").Value = '2026-02-09 07:01:52'
$ws.Range("B9").Value = '【急募】多言語動画吹替・字幕一括生成システム開発'
$ws.Range("C9").Value = 'システム開発'
$ws.Range("D9").Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Range("E9").Value = '期限情報なし'
$ws.Range("F9").Value = 'https://www.lancers.jp/work/detail/5488276'
$ws.Range("G9").Value = 118
$ws.Range("H9").Value = '◆開発,システム開発'

# Row 10
$ws.Range("A10").Value = '2026-02-09 07:01:52'
$ws.Range("B10").Value = '養鰻管理Excelの判断ロジック(給餌)を理解し、継続的に伴走できる方を募集'
$ws.Range("C10").Value = 'システム開発'
$ws.Range("D10").Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Range("E10").Value = '期限情報なし'
$ws.Range("F10").Value = 'https://www.lancers.jp/work/detail/5488109'
$ws.Range("G10").Value = 38
$ws.Range("H10").Value = '◇管理'

# Re-create hyperlinks for F2:F10 and restore the Hyperlink style
$ws.Hyperlinks.Add($ws.Range("F2"), 'https://www.lancers.jp/work/detail/5488301')
$ws.Range("F2").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F3"), 'https://www.lancers.jp/work/detail/5488299')
$ws.Range("F3").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F4"), 'https://www.lancers.jp/work/detail/5488286')
$ws.Range("F4").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F5"), 'https://www.lancers.jp/work/detail/5488266')
$ws.Range("F5").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F6"), 'https://www.lancers.jp/work/detail/5488392')
$ws.Range("F6").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F7"), 'https://www.lancers.jp/work/detail/5488271')
$ws.Range("F7").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F8"), 'https://www.lancers.jp/work/detail/5488168')
$ws.Range("F8").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F9"), 'https://www.lancers.jp/work/detail/5488276')
$ws.Range("F9").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F10"), 'https://www.lancers.jp/work/detail/5488109')
$ws.Range("F10").Style = "Hyperlink"

# Column width adjustments (ColumnWidth API adds ~0.8333 padding vs. stored width)
$ws.Columns.Item(4).ColumnWidth = 28 - 5/6
$ws.Columns.Item(8).ColumnWidth = 19 - 5/6
